$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "37.848.24"
Set-TextValue $ws "E2" "  -0.24%  "
Set-TextValue $ws "D3" "2.084.36"
Set-TextValue $ws "E3" "  -0.47%  "
Set-TextValue $ws "E4" "  +0.01%  "
Set-TextValue $ws "D5" "233.56"
Set-TextValue $ws "E5" "  +0.10%  "
Set-TextValue $ws "E6" "  -0.24%  "
Set-TextValue $ws "E7" "  +2.80%  "
Set-TextValue $ws "E8" "  -0.07%  "
Set-TextValue $ws "E9" "  +1.50%  "
Set-TextValue $ws "D10" "0.0789"
Set-TextValue $ws "E10" "  +0.69%  "
Set-TextValue $ws "E11" "  +1.77%  "
Set-TextValue $ws "D12" "2.391.99"
Set-TextValue $ws "E12" "  -0.25%  "
Set-TextValue $ws "D13" "14.77"
Set-TextValue $ws "E13" "  +2.11%  "
Set-TextValue $ws "E14" "  +0.15%  "
Set-TextValue $ws "E15" "  +1.03%  "
Set-TextValue $ws "D16" "5.34"
Set-TextValue $ws "E16" "  +1.52%  "
Set-TextValue $ws "D17" "2.107.92"
Set-TextValue $ws "E17" "  +0.49%  "
Set-TextValue $ws "D18" "37.780.71"
Set-TextValue $ws "E18" "  -0.32%  "
Set-TextValue $ws "D19" "6.16"
Set-TextValue $ws "E19" "  +0.43%  "
Set-TextValue $ws "D20" "71.64"
Set-TextValue $ws "E20" "  +1.01%  "
Set-TextValue $ws "D21" "0.0₃0849"
Set-TextValue $ws "E21" "  +3.19%  "
Set-TextValue $ws "D22" "228.01"
Set-TextValue $ws "E23" "  -0.07%  "
Set-TextValue $ws "E24" "  -0.74%  "
Set-TextValue $ws "E25" "  +0.29%  "
Set-TextValue $ws "D26" "171.77"
Set-TextValue $ws "E26" "  +0.64%  "
Set-TextValue $ws "D27" "9.29"
Set-TextValue $ws "E27" "  +3.68%  "
Set-TextValue $ws "E28" "  -2.55%  "
Set-TextValue $ws "E29" "  -1.33%  "
Set-TextValue $ws "D30" "19.51"
Set-TextValue $ws "E30" "  +0.02%  "
Set-TextValue $ws "E31" "  +1.71%  "
Set-TextValue $ws "D33" "0.0634"
Set-TextValue $ws "E33" "  +0.53%  "
Set-TextValue $ws "D34" "4.70"
Set-TextValue $ws "E34" "  +2.21%  "
Set-TextValue $ws "E35" "  -0.87%  "
Set-TextValue $ws "E36" "  +0.38%  "
Set-TextValue $ws "E37" "  -0.85%  "
Set-TextValue $ws "E38" "  +0.00%  "
Set-TextValue $ws "E39" "  -0.78%  "
Set-TextValue $ws "D40" "0.0985"
Set-TextValue $ws "E40" "  -1.69%  "
Set-TextValue $ws "D41" "99.47"
Set-TextValue $ws "E41" "  +2.10%  "
Set-TextValue $ws "D43" "17.10"
Set-TextValue $ws "E43" "  +8.22%  "
Set-TextValue $ws "D44" "2.90"
Set-TextValue $ws "E44" "  -1.25%  "
Set-TextValue $ws "D45" "1.450.78"
Set-TextValue $ws "E45" "  -0.30%  "
Set-TextValue $ws "E46" "  -1.19%  "
Set-TextValue $ws "D47" "4.17"
Set-TextValue $ws "E47" "  +3.51%  "
Set-TextValue $ws "E48" "  +0.56%  "
Set-TextValue $ws "E49" "  -0.43%  "
Set-TextValue $ws "D50" "3.00"
Set-TextValue $ws "D51" "2.277.28"
Set-TextValue $ws "E51" "  -0.49%  "
